$d = $word.ActiveDocument

function Remove-ParagraphByText($text) {
    $rng = $d.Content
    $found = $rng.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return
    }
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Start -le $rng.Start -and $p.Range.End -ge $rng.End) {
            $p.Range.Delete()
            return
        }
    }
}

# 1. Remove the "_GoBack" bookmark left over at the top of the document.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Remove the "Save to temporary file before renaming to desired file" bullet.
Remove-ParagraphByText("Save to temporary file before renaming to desired file")

# 3. Remove the two "Nice to haves" sub-bullets about table resizing / word wrapping.
Remove-ParagraphByText("Table column & row manual resizing is smooth & shown real-time")
Remove-ParagraphByText("Table cell improved text word wrapping (as in Qt)")
